$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells M1:P1
$ws.Range("M1").Value = "%Change_Test"
$ws.Range("N1").Value = "%Change_Control"
$ws.Range("O1").Value = "%Change_Diff"
$ws.Range("P1").Value = "Direction"

# Copy the header style from an existing header cell (A1) onto the new headers
$ws.Range("A1").Copy()
$ws.Range("M1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2
$ws.Range("M2").Value = -30.70934256055362
$ws.Range("N2").Value = -1.556482670089866
$ws.Range("O2").Value = -29.15285989046376
$ws.Range("P2").Value = "↓ Worse"

# Row 3
$ws.Range("M3").Value = 1.263547967502731
$ws.Range("N3").Value = -1.106314900155447
$ws.Range("O3").Value = 2.369862867658178
$ws.Range("P3").Value = "↑ Better"

# Row 4
$ws.Range("M4").Value = -17.05933247265292
$ws.Range("N4").Value = -16.00011230023173
$ws.Range("O4").Value = -1.059220172421181
$ws.Range("P4").Value = "↓ Worse"

# Row 5
$ws.Range("M5").Value = -0.06470382003188843
$ws.Range("N5").Value = -0.06041741005078369
$ws.Range("O5").Value = -0.004286409981104737
$ws.Range("P5").Value = "↓ Worse"

# Row 6
$ws.Range("M6").Value = 23.30296209323099
$ws.Range("N6").Value = -21.91155812647189
$ws.Range("O6").Value = 45.21452021970288
$ws.Range("P6").Value = "↑ Better"

# Row 7
$ws.Range("M7").Value = -28.67067849865575
$ws.Range("N7").Value = -29.84685902058545
$ws.Range("O7").Value = 1.176180521929695
$ws.Range("P7").Value = "↑ Better"

# Row 8
$ws.Range("M8").Value = 90.04344012309467
$ws.Range("N8").Value = 147.4293694420277
$ws.Range("O8").Value = -57.38592931893304
$ws.Range("P8").Value = "↓ Worse"
